$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "42.843.00"
$ws.Range("E2").Value2 = "  +0.97%  "
$ws.Range("D3").Value2 = "2.530.58"
$ws.Range("E3").Value2 = "  +0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("E4").Value2 = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "316.95"
$ws.Range("E5").Value2 = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "96.23"
$ws.Range("E6").Value2 = "  +1.51%  "
$ws.Range("E7").Value2 = "  -0.17%  "
$ws.Range("E8").Value2 = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.538"
$ws.Range("E9").Value2 = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "35.50"
$ws.Range("E10").Value2 = "  -1.97%  "
$ws.Range("E11").Value2 = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "7.48"
$ws.Range("E12").Value2 = "  -2.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.108"
$ws.Range("E13").Value2 = "  -4.88%  "
$ws.Range("D14").Value2 = "2.917.95"
$ws.Range("E14").Value2 = "  +0.15%  "
$ws.Range("D15").Value2 = "2.522.24"
$ws.Range("E15").Value2 = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "15.00"
$ws.Range("E16").Value2 = "  -4.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.846"
$ws.Range("E17").Value2 = "  -2.00%  "
$ws.Range("D18").Value2 = "42.857.59"
$ws.Range("E18").Value2 = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.81"
$ws.Range("E19").Value2 = "  +2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "12.50"
$ws.Range("E20").Value2 = "  -3.86%  "
$ws.Range("D21").Value2 = "0.0₃0961"
$ws.Range("E21").Value2 = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "69.45"
$ws.Range("E22").Value2 = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "253.35"
$ws.Range("E23").Value2 = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.95"
$ws.Range("E24").Value2 = "  -0.36%  "
$ws.Range("E25").Value2 = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "26.68"
$ws.Range("E26").Value2 = "  -0.65%  "
$ws.Range("E27").Value2 = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "40.44"
$ws.Range("E29").Value2 = "  +2.53%  "
$ws.Range("E30").Value2 = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "5.87"
$ws.Range("E31").Value2 = "  -1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "156.62"
$ws.Range("E32").Value2 = "  +0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "2.71"
$ws.Range("E33").Value2 = "  +3.46%  "
$ws.Range("B34").Value2 = "Celestia"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "19.12"
$ws.Range("E34").Value2 = "  +1.49%  "
$ws.Range("B35").Value2 = "LidoDAOToken"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.32"
$ws.Range("E35").Value2 = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.0794"
$ws.Range("E36").Value2 = "  +1.15%  "
$ws.Range("E37").Value2 = "  -2.30%  "
$ws.Range("E38").Value2 = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.44"
$ws.Range("E39").Value2 = "  +4.85%  "
$ws.Range("E40").Value2 = "  -0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "21.75"
$ws.Range("E41").Value2 = "  -8.60%  "
$ws.Range("B42").Value2 = "RenderToken"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "3.79"
$ws.Range("E42").Value2 = "  -0.84%  "
$ws.Range("B43").Value2 = "FirstDigitalUSD"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "1.00"
$ws.Range("E43").Value2 = "  -0.29%  "
$ws.Range("E44").Value2 = "  +1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "3.26"
$ws.Range("E45").Value2 = "  -2.28%  "
$ws.Range("D46").Value2 = "1.990.57"
$ws.Range("E46").Value2 = "  -1.77%  "
$ws.Range("E47").Value2 = "  +2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "84.47"
$ws.Range("E48").Value2 = "  +0.11%  "
$ws.Range("D49").Value2 = "2.772.35"
$ws.Range("E49").Value2 = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "74.31"
$ws.Range("E50").Value2 = "  +1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "104.19"
$ws.Range("E51").Value2 = "  +2.63%  "
